# Finalize the data curation:
#  - rename the placeholder row labels (1..11) in column A to the real
#    variable names now that the dataset has been curated
#  - update the measured values in columns B (SampleSize/English col),
#    C and D to the corrected figures
#  - drop the final sample row (previously row 12 / label "11"), which
#    didn't survive curation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("1" -> Acidity)
$ws.Range("A2").Value = "Acidity"
$ws.Range("B2").Value = 14.0
$ws.Range("C2").Value = 10.0
$ws.Range("D2").Value = 28.0

# Row 3 ("2" -> DeltaAcidity)
$ws.Range("A3").Value = "DeltaAcidity"
$ws.Range("B3").Value = 14.0
$ws.Range("C3").Value = 8.0
$ws.Range("D3").Value = 20.0

# Row 4 ("3" -> SSC)
$ws.Range("A4").Value = "SSC"
$ws.Range("B4").Value = 15.0
$ws.Range("C4").Value = 10.0
$ws.Range("D4").Value = 29.0

# Row 5 ("4" -> Firmness)
$ws.Range("A5").Value = "Firmness"
$ws.Range("B5").Value = 15.0
$ws.Range("C5").Value = 10.0
$ws.Range("D5").Value = 29.0

# Row 6 ("5" -> Weight)
$ws.Range("A6").Value = "Weight"
$ws.Range("B6").Value = 15.0
$ws.Range("C6").Value = 10.0
$ws.Range("D6").Value = 29.0

# Row 7 ("6" -> Juiciness)
$ws.Range("A7").Value = "Juiciness"
$ws.Range("B7").Value = 12.0
$ws.Range("C7").Value = 5.0
$ws.Range("D7").Value = 13.0

# Row 8 ("7" -> PhenolicContent)
$ws.Range("A8").Value = "PhenolicContent"
$ws.Range("B8").Value = 12.0
$ws.Range("C8").Value = 3.0
$ws.Range("D8").Value = 12.0

# Row 9 ("8" -> HarvestDate)
$ws.Range("A9").Value = "HarvestDate"
$ws.Range("B9").Value = 15.0
$ws.Range("C9").Value = 10.0
$ws.Range("D9").Value = 29.0

# Row 10 ("9" -> FloweringDate)
$ws.Range("A10").Value = "FloweringDate"
$ws.Range("B10").Value = 15.0
$ws.Range("C10").Value = 10.0
$ws.Range("D10").Value = 29.0

# Row 11 ("10" -> Softening)
$ws.Range("A11").Value = "Softening"
$ws.Range("B11").Value = 12.0
$ws.Range("C11").Value = 6.0
$ws.Range("D11").Value = 23.0

# Row 12 ("11") is dropped entirely from the curated dataset
$ws.Rows(12).Delete()
